# Apply updated sensor readings (columns C:H) for the "falling" dataset sheet,
# shifting all numeric data down by two rows' worth of new samples and
# dropping the final row so the used range shrinks from A1:H22 to A1:H21.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New values for each row (columns C,D,E,F,G,H), rows 2-21 in final layout.
$data = @{
    2  = @(-1.921436786651612, 1.638089656829834, 0.3902863562107086, -0.1360424668951469, 0.3186911859295585, 0.05561650341207323)
    3  = @(-1.84520435333252, 1.668948650360107, 0.1994338035583496, -0.2708076590841467, 0.3116106986999512, 0.0008329986171288634)
    4  = @(-1.829729080200196, 1.47260046005249, 0.5425764322280884, -0.3139847218990326, 0.365755558013916, -0.0809396430850029)
    5  = @(-1.84237813949585, 1.411740303039551, 1.237700462341309, -0.09199075129899148, 0.302864204753529, -0.2070000191981139)
    6  = @(-2.159674644470215, 1.136321544647217, 3.083052396774292, 0.01230060461569882, -0.3993948210369431, -0.9949190426956491)
    7  = @(-1.937846660614014, 1.245656490325928, 3.884034872055054, 0.3243139045820992, -0.5972596244378523, -1.587486072020097)
    8  = @(-3.061091899871826, 0.4542388916015625, 2.74072003364563, -0.9380390860817602, -1.513376929543232, -3.357259793715038)
    9  = @(-2.979402542114258, -0.1682605743408203, 2.131577730178833, -0.2036124901338088, -0.2562440525401761, -0.1245054331692872)
    10 = @(-3.620201349258423, -1.011659622192383, 2.961694955825806, 0.2954227382486529, 0.008163343776357967, 0.4832359877499712)
    11 = @(-5.750590324401856, -1.630428791046143, 4.153227806091309, -0.3702815256335535, -0.7731332995674872, -1.987241875041621)
    12 = @(-4.631756782531738, 0.1280508041381836, 0.0708565711975097, -0.1218537308953028, -0.1149676279588187, -0.5308419032530356)
    13 = @(-6.701959609985352, 13.03308868408203, -3.955403804779053, -0.7092561884359876, -1.314401550726458, -0.9030254652554355)
    14 = @(-9.679704666137695, 1.827943325042725, 0.2789157629013061, -0.9390387182885965, -0.6636356765573627, 0.2282692268490727)
    15 = @(3.086193084716797, -2.266827583312988, 2.888416528701782, -0.3692680299282074, -2.486374378204346, 2.66978645324707)
    16 = @(1.258205413818359, 1.049188137054443, -6.970268726348877, 0.7195992388508516, 0.8889339186928424, 0.3780700483105441)
    17 = @(-2.633459091186523, 5.726268291473389, 6.590025901794434, -1.744853258132938, -2.628608725287703, 2.057380207560282)
    18 = @(4.84531307220459, -0.8891797065734863, -5.761303424835205, -1.074428666721691, -0.7690515951676815, 1.639270782470704)
    19 = @(0.9592771530151368, 3.222768306732178, 1.57231593132019, -0.2539533461359424, -0.02276861396702762, 0.6248873986981154)
    20 = @(1.009872436523438, -0.1888983249664306, 0.460105299949646, -0.05406157367608758, 0.4298825345256099, -0.1093171049248081)
    21 = @(-3.207517623901367, 0.4280340671539306, 0.7658971548080444, 0.0573519143191249, 0.5145012542012748, -0.1331963850693274)
}

$cols = @("C", "D", "E", "F", "G", "H")
foreach ($row in $data.Keys) {
    $vals = $data[$row]
    for ($i = 0; $i -lt $cols.Length; $i++) {
        $ws.Range("$($cols[$i])$row").Value = $vals[$i]
    }
}

# Remove the now-obsolete last row (old row 22) so the sheet's used range
# shrinks to A1:H21, matching the reduced dataset.
$ws.Rows.Item(22).Delete()
